$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$null = $ws.Range("Q4:Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Font.Bold = $true
